# "Season up to 1/17"
#
# The NBA game that was previously the next scheduled matchup (UTA on
# 45306, the first row of the "Next" sheet) has now been played, so it
# moves from the "Next" sheet into the "Games" sheet as a new completed
# game (row 41), and is removed from the "Next" sheet (all remaining
# rows shift up by one).

$wb = $excel.ActiveWorkbook

$gamesWs = $wb.Worksheets.Item("Games")
$nextWs  = $wb.Worksheets.Item("Next")

# --- 1) Append the newly completed game as row 41 on the "Games" sheet ---
$newRow = 41

$gamesWs.Cells.Item($newRow, 1).Value  = 40       # Game
$gamesWs.Cells.Item($newRow, 2).Value  = 45306     # Date
$gamesWs.Cells.Item($newRow, 2).NumberFormat = $gamesWs.Cells.Item($newRow - 1, 2).NumberFormat
$gamesWs.Cells.Item($newRow, 3).Value  = -2        # Streak
$gamesWs.Cells.Item($newRow, 4).Value  = 105       # Pts
$gamesWs.Cells.Item($newRow, 5).Value  = 97.8      # Pace
$gamesWs.Cells.Item($newRow, 6).Value  = 0.447     # eFG
$gamesWs.Cells.Item($newRow, 7).Value  = 10.9      # TOV
$gamesWs.Cells.Item($newRow, 8).Value  = 33.3      # ORB
$gamesWs.Cells.Item($newRow, 9).Value  = 0.211     # FTR
$gamesWs.Cells.Item($newRow, 10).Value = 107.3     # ORT
$gamesWs.Cells.Item($newRow, 11).Value = "UTA"     # OppID
$gamesWs.Cells.Item($newRow, 12).Value = 132       # OppPts
$gamesWs.Cells.Item($newRow, 13).Value = 0.623     # OppeFG
$gamesWs.Cells.Item($newRow, 14).Value = 12.7      # OppTOV
$gamesWs.Cells.Item($newRow, 15).Value = 28.6      # OppORB
$gamesWs.Cells.Item($newRow, 16).Value = 0.383     # OppFTR
$gamesWs.Cells.Item($newRow, 17).Value = 134.9     # OppORT
$gamesWs.Cells.Item($newRow, 18).Value = 0         # Location
$gamesWs.Cells.Item($newRow, 19).Value = 0         # Target

# --- 2) Remove that game from the "Next" sheet; remaining rows shift up ---
$nextWs.Rows.Item(2).Delete()
